$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" column (I) was added to the table that currently runs
# through column H (years 2016-2020). Copy the number formatting/styles
# from column H (rows 4-25, the header + all data rows) into the new
# column I so each new cell inherits the same look as its neighbour.
$ws.Range("H4:H25").Copy()
$ws.Range("I4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header
$ws.Range("I4").Value = 2021

# Data rows (row => 2021 value). Rows 6, 9, 13 and 16 are section/group
# header rows with no numeric value, so column I stays blank there too
# (formatting only, already handled by the PasteSpecial above).
$values = [ordered]@{
  5  = 48.5
  7  = 48.8
  8  = 48.2
  10 = 58.2
  11 = 42.4
  12 = 40.7
  14 = 41.5
  15 = 52.6
  17 = 67.1
  18 = 62
  19 = 46.9
  20 = 55.8
  21 = 42.7
  22 = 48.3
  23 = 39.7
  24 = 38.1
  25 = 44.7
}

foreach ($row in $values.Keys) {
  $ws.Range("I$row").Value = $values[$row]
}

# The author's selection had moved off the old A14:C15 block back to the
# default top-left cell by the time this was saved.
$ws.Range("A1").Select()
